$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 10489
$ws.Range("F3").Value = 245
$ws.Range("F4").Value = 1160
$ws.Range("F5").Value = 1033
$ws.Range("F6").Value = 820
$ws.Range("F7").Value = 264
$ws.Range("F9").Value = 336
$ws.Range("F10").Value = 1095
$ws.Range("F12").Value = 145
$ws.Range("F13").Value = 846
$ws.Range("F14").Value = 366
$ws.Range("F15").Value = 1779
$ws.Range("F17").Value = 894
$ws.Range("F18").Value = 799
$ws.Range("F19").Value = 529
$ws.Range("F20").Value = 766
$ws.Range("F21").Value = 862
$ws.Range("F22").Value = 42
$ws.Range("F23").Value = 265
$ws.Range("F24").Value = 84
$ws.Range("F25").Value = 593
$ws.Range("F26").Value = 606
$ws.Range("F27").Value = 101
$ws.Range("F28").Value = 314
$ws.Range("F29").Value = 992
$ws.Range("F31").Value = 477
$ws.Range("F32").Value = 144
$ws.Range("F33").Value = 243
$ws.Range("F34").Value = 209
$ws.Range("F35").Value = 533
$ws.Range("F36").Value = 1620
$ws.Range("F37").Value = 359
$ws.Range("F38").Value = 46
$ws.Range("F39").Value = 1381
$ws.Range("F40").Value = 399
$ws.Range("F43").Value = 75
$ws.Range("F44").Value = 40
$ws.Range("F45").Value = 65
$ws.Range("F46").Value = 64
$ws.Range("F47").Value = 32
$ws.Range("F48").Value = 6
$ws.Range("F49").Value = 75

$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 179
$ws.Range("F6").Value = 12
$ws.Range("F7").Value = 71
$ws.Range("F13").Value = 101
$ws.Range("F16").Value = 9
$ws.Range("F17").Value = 9

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 2142
$ws.Range("F3").Value = 599
$ws.Range("F4").Value = 518

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2142
$ws.Range("F3").Value = 599
$ws.Range("F5").Value = 10489
$ws.Range("F6").Value = 245
$ws.Range("F7").Value = 1160
$ws.Range("F8").Value = 518
$ws.Range("F9").Value = 1033
$ws.Range("F10").Value = 820
$ws.Range("F11").Value = 179
$ws.Range("F12").Value = 336
$ws.Range("F13").Value = 1095
$ws.Range("F14").Value = 12
$ws.Range("F15").Value = 145
$ws.Range("F16").Value = 846
$ws.Range("F17").Value = 366
$ws.Range("F18").Value = 1779
$ws.Range("F20").Value = 894
$ws.Range("F21").Value = 799
$ws.Range("F22").Value = 529
$ws.Range("F23").Value = 766
$ws.Range("F24").Value = 862
$ws.Range("F25").Value = 42
$ws.Range("F26").Value = 265
$ws.Range("F27").Value = 84
$ws.Range("F28").Value = 593
$ws.Range("F31").Value = 606
$ws.Range("F32").Value = 101
$ws.Range("F33").Value = 314
$ws.Range("F34").Value = 992
$ws.Range("F37").Value = 477
$ws.Range("F38").Value = 144
$ws.Range("F40").Value = 359
$ws.Range("F41").Value = 1381
$ws.Range("F42").Value = 399
$ws.Range("F46").Value = 75
$ws.Range("F47").Value = 65
$ws.Range("F48").Value = 75
$ws.Range("F49").Value = 9
